$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Style = "Normal"

$metaText = "Meta description: Experience classic gameplay with a modern twist in 40 Joker Staxx: 40 lines. Read our review and play for free today."
$metaRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End - 1)
$metaRange.Text = $metaText

$labelLen = [int]("Meta description".Length)
$boldRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + $labelLen)
$boldRange.Bold = 1

# ------------------------------------------------------------------
# 2) Remove the duplicated bold "Play 40 Joker Staxx..." paragraph
#    that used to sit right before the italic meta-description /
#    image-prompt paragraph near the end of the document, and update
#    that italic paragraph's text to the new image-generation prompt.
# ------------------------------------------------------------------
$titleText = "Play 40 Joker Staxx: 40 lines for Free - Read Our Review"

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text.Trim()
    if ($i -gt 2 -and $txt -eq $titleText) {
        $para.Range.Delete()
        break
    }
}

$imagePrompt = "Create a feature image for `"40 Joker Staxx: 40 Lines`" that highlights the game's modern twist on retro themes. The image should be in a cartoon style, featuring a happy Maya warrior wearing glasses to represent the game's simple yet fun gameplay. The warrior could be holding a classic fruit or a gold ingot, two symbols that represent the game's payout potential. The background could be a mix of retro and modern elements, such as neon lights and classic arcade machines. The overall tone of the image should be vibrant and lively, reflecting the excitement of playing the game."

$oldDescription = "Experience classic gameplay with a modern twist in 40 Joker Staxx: 40 lines. Read our review and play for free today."
$oldDescriptionTrimmed = $oldDescription.Trim()

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text.Trim()
    if ($txt -eq $oldDescriptionTrimmed) {
        $target = $d.Range($para.Range.Start, $para.Range.End - 1)
        $target.Text = $imagePrompt
        break
    }
}
